# Auto-generated Excel COM-interop script
# Applies numeric value updates (and a few cell add/remove operations)
# to the "Ultros_Profits" workbook market-data tables across all 8 job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ALC_edits = @(
  @{ Cell = "H28"; Value = 2102.4546 },
  @{ Cell = "I28"; Value = 2112.7 },
  @{ Cell = "K28"; Value = 2112.7 },
  @{ Cell = "M28"; Value = -1627.7 },
  @{ Cell = "H40"; Value = 6125.143 },
  @{ Cell = "J40"; Value = 7857 },
  @{ Cell = "L40"; Value = 7857 },
  @{ Cell = "N40"; Value = -8207 },
  @{ Cell = "H41"; Value = 890.86206 },
  @{ Cell = "I41"; Value = 824.6667 },
  @{ Cell = "J41"; Value = 1064.625 },
  @{ Cell = "K41"; Value = 824.6667 },
  @{ Cell = "L41"; Value = 1064.625 },
  @{ Cell = "M41"; Value = -384.6667 },
  @{ Cell = "N41"; Value = -1944.625 },
  @{ Cell = "H42"; Value = 200175.39 },
  @{ Cell = "I42"; Value = 260184.3 },
  @{ Cell = "J42"; Value = 145.66667 },
  @{ Cell = "K42"; Value = 780552.8999999999 },
  @{ Cell = "L42"; Value = 437.00001 },
  @{ Cell = "M42"; Value = -780322.8999999999 },
  @{ Cell = "N42"; Value = -897.00001 },
  @{ Cell = "H62"; Value = 55558056 },
  @{ Cell = "I62"; Value = 55558056 },
  @{ Cell = "K62"; Value = 55558056 },
  @{ Cell = "M62"; Value = -55557432 },
  @{ Cell = "H65"; Value = 55558056 },
  @{ Cell = "I65"; Value = 55558056 },
  @{ Cell = "K65"; Value = 277790280 },
  @{ Cell = "M65"; Value = -277787160 },
  @{ Cell = "H87"; Value = 19999.715 },
  @{ Cell = "J87"; Value = 19999.715 },
  @{ Cell = "L87"; Value = 19999.715 },
  @{ Cell = "N87"; Value = -22495.715 },
  @{ Cell = "H88"; Value = 2536.1052 },
  @{ Cell = "I88"; Value = 4299.778 },
  @{ Cell = "J88"; Value = 948.8 },
  @{ Cell = "K88"; Value = 4299.778 },
  @{ Cell = "L88"; Value = 948.8 },
  @{ Cell = "M88"; Value = -3893.778 },
  @{ Cell = "N88"; Value = -1760.8 },
  @{ Cell = "H90"; Value = 19999.715 },
  @{ Cell = "J90"; Value = 19999.715 },
  @{ Cell = "L90"; Value = 59999.145 },
  @{ Cell = "N90"; Value = -72479.145 },
  @{ Cell = "H91"; Value = 2536.1052 },
  @{ Cell = "I91"; Value = 4299.778 },
  @{ Cell = "J91"; Value = 948.8 },
  @{ Cell = "K91"; Value = 4299.778 },
  @{ Cell = "L91"; Value = 948.8 },
  @{ Cell = "M91"; Value = -2895.778 },
  @{ Cell = "N91"; Value = -3756.8 },
  @{ Cell = "H107"; Value = 2325.8 },
  @{ Cell = "J107"; Value = 267.5 },
  @{ Cell = "L107"; Value = 267.5 },
  @{ Cell = "N107"; Value = -4107.5 },
  @{ Cell = "H138"; Value = 2874.38 },
  @{ Cell = "I138"; Value = 1859.3334 },
  @{ Cell = "K138"; Value = 5578.0002 },
  @{ Cell = "M138"; Value = -438.0002000000004 }
)
foreach ($e in $ALC_edits) {
  $ws.Range($e.Cell).Value = $e.Value
}

$ws = $wb.Worksheets.Item("ARM")

$ARM_edits = @(
  @{ Cell = "H32"; Value = 47654012 },
  @{ Cell = "I32"; Value = 66697148 },
  @{ Cell = "J32"; Value = 46165.5 },
  @{ Cell = "K32"; Value = 66697148 },
  @{ Cell = "L32"; Value = 46165.5 },
  @{ Cell = "M32"; Value = -66696861 },
  @{ Cell = "N32"; Value = -46739.5 },
  @{ Cell = "H74"; Value = 1867.2106 },
  @{ Cell = "I74"; Value = 1379.2667 },
  @{ Cell = "K74"; Value = 1379.2667 },
  @{ Cell = "M74"; Value = -505.2666999999999 },
  @{ Cell = "H77"; Value = 1867.2106 },
  @{ Cell = "I77"; Value = 1379.2667 },
  @{ Cell = "K77"; Value = 6896.3335 },
  @{ Cell = "M77"; Value = -2528.3335 },
  @{ Cell = "H88"; Value = 3975016.2 },
  @{ Cell = "I88"; Value = 11294.454 },
  @{ Cell = "K88"; Value = 11294.454 },
  @{ Cell = "M88"; Value = -10888.454 },
  @{ Cell = "H91"; Value = 3975016.2 },
  @{ Cell = "I91"; Value = 11294.454 },
  @{ Cell = "K91"; Value = 11294.454 },
  @{ Cell = "M91"; Value = -9890.454 },
  @{ Cell = "H102"; Value = 20834468 },
  @{ Cell = "I102"; Value = 1208.5333 },
  @{ Cell = "K102"; Value = 1208.5333 },
  @{ Cell = "M102"; Value = 413.4666999999999 },
  @{ Cell = "H132"; Value = 4656.5264 },
  @{ Cell = "I132"; Value = 3807.3076 },
  @{ Cell = "K132"; Value = 11421.9228 },
  @{ Cell = "M132"; Value = -8891.9228 }
)
foreach ($e in $ARM_edits) {
  $ws.Range($e.Cell).Value = $e.Value
}

$ws = $wb.Worksheets.Item("BSM")

$BSM_edits = @(
  @{ Cell = "H105"; Value = 3429.75 },
  @{ Cell = "I105"; Value = 1901.6666 },
  @{ Cell = "J105"; Value = 4957.8335 },
  @{ Cell = "K105"; Value = 1901.6666 },
  @{ Cell = "L105"; Value = 4957.8335 },
  @{ Cell = "M105"; Value = -154.6666 },
  @{ Cell = "N105"; Value = -8451.8335 }
)
foreach ($e in $BSM_edits) {
  $ws.Range($e.Cell).Value = $e.Value
}

$ws = $wb.Worksheets.Item("CRP")

$CRP_edits = @(
  @{ Cell = "H31"; Value = 2630.4707 },
  @{ Cell = "I31"; Value = 1992.1364 },
  @{ Cell = "K31"; Value = 1992.1364 },
  @{ Cell = "M31"; Value = -1697.1364 },
  @{ Cell = "H34"; Value = 2630.4707 },
  @{ Cell = "I34"; Value = 1992.1364 },
  @{ Cell = "K34"; Value = 1992.1364 },
  @{ Cell = "M34"; Value = -1790.1364 },
  @{ Cell = "H93"; Value = 15268.75 },
  @{ Cell = "I93"; Value = 3164.2856 },
  @{ Cell = "K93"; Value = 3164.2856 },
  @{ Cell = "M93"; Value = -1292.2856 }
)
foreach ($e in $CRP_edits) {
  $ws.Range($e.Cell).Value = $e.Value
}

$ws = $wb.Worksheets.Item("CUL")

$CUL_edits = @(
  @{ Cell = "H37"; Value = 211956 },
  @{ Cell = "J37"; Value = 211956 },
  @{ Cell = "L37"; Value = 635868 },
  @{ Cell = "N37"; Value = -636092 },
  @{ Cell = "H87"; Value = 25003852 },
  @{ Cell = "I87"; Value = 25003852 },
  @{ Cell = "J87"; Value = 0 },
  @{ Cell = "K87"; Value = 75011556 },
  @{ Cell = "L87"; Value = 0 },
  @{ Cell = "M87"; Value = -75010308 },
  @{ Cell = "H90"; Value = 25003852 },
  @{ Cell = "I90"; Value = 25003852 },
  @{ Cell = "J90"; Value = 0 },
  @{ Cell = "K90"; Value = 225034668 },
  @{ Cell = "L90"; Value = 0 },
  @{ Cell = "M90"; Value = -225028428 },
  @{ Cell = "H114"; Value = 1371.4667 },
  @{ Cell = "I114"; Value = 497.5 },
  @{ Cell = "J114"; Value = 1689.2727 },
  @{ Cell = "K114"; Value = 1492.5 },
  @{ Cell = "L114"; Value = 5067.8181 },
  @{ Cell = "M114"; Value = 1761.5 },
  @{ Cell = "N114"; Value = -11575.8181 },
  @{ Cell = "H117"; Value = 988.4167 },
  @{ Cell = "I117"; Value = 531.5 },
  @{ Cell = "K117"; Value = 1594.5 },
  @{ Cell = "M117"; Value = 1847.5 },
  @{ Cell = "H122"; Value = 2004.8235 },
  @{ Cell = "I122"; Value = 1733.1666 },
  @{ Cell = "J122"; Value = 2153 },
  @{ Cell = "K122"; Value = 15598.4994 },
  @{ Cell = "L122"; Value = 19377 },
  @{ Cell = "M122"; Value = -13148.4994 },
  @{ Cell = "N122"; Value = -24277 },
  @{ Cell = "H129"; Value = 1532.5 },
  @{ Cell = "I129"; Value = 988.8889 },
  @{ Cell = "J129"; Value = 2231.4285 },
  @{ Cell = "K129"; Value = 2966.6667 },
  @{ Cell = "L129"; Value = 6694.2855 },
  @{ Cell = "M129"; Value = 2033.3333 },
  @{ Cell = "N129"; Value = -16694.2855 }
)
foreach ($e in $CUL_edits) {
  $ws.Range($e.Cell).Value = $e.Value
}

$CUL_deletes = @("N87","N90")
foreach ($c in $CUL_deletes) {
  $ws.Range($c).ClearContents()
}

$ws = $wb.Worksheets.Item("GSM")

$GSM_edits = @(
  @{ Cell = "H123"; Value = 36833.332 },
  @{ Cell = "J123"; Value = 36833.332 },
  @{ Cell = "L123"; Value = 36833.332 },
  @{ Cell = "N123"; Value = -41733.332 },
  @{ Cell = "H132"; Value = 4482.1 },
  @{ Cell = "I132"; Value = 4130.5264 },
  @{ Cell = "K132"; Value = 12391.5792 },
  @{ Cell = "M132"; Value = -9861.5792 }
)
foreach ($e in $GSM_edits) {
  $ws.Range($e.Cell).Value = $e.Value
}

$ws = $wb.Worksheets.Item("LTW")

$LTW_edits = @(
  @{ Cell = "H61"; Value = 3471.5356 },
  @{ Cell = "I61"; Value = 1923.8235 },
  @{ Cell = "J61"; Value = 5863.4546 },
  @{ Cell = "K61"; Value = 1923.8235 },
  @{ Cell = "L61"; Value = 5863.4546 },
  @{ Cell = "M61"; Value = -1721.8235 },
  @{ Cell = "N61"; Value = -6267.4546 },
  @{ Cell = "H96"; Value = 0 },
  @{ Cell = "J96"; Value = 0 },
  @{ Cell = "L96"; Value = 0 },
  @{ Cell = "H104"; Value = 18285.715 },
  @{ Cell = "J104"; Value = 18285.715 },
  @{ Cell = "L104"; Value = 18285.715 },
  @{ Cell = "N104"; Value = -25273.715 },
  @{ Cell = "H106"; Value = 25000 },
  @{ Cell = "J106"; Value = 25000 },
  @{ Cell = "L106"; Value = 25000 },
  @{ Cell = "N106"; Value = -27524 },
  @{ Cell = "H113"; Value = 3471.5356 },
  @{ Cell = "I113"; Value = 1923.8235 },
  @{ Cell = "J113"; Value = 5863.4546 },
  @{ Cell = "K113"; Value = 1923.8235 },
  @{ Cell = "L113"; Value = 5863.4546 },
  @{ Cell = "M113"; Value = 246.1765 },
  @{ Cell = "N113"; Value = -10203.4546 }
)
foreach ($e in $LTW_edits) {
  $ws.Range($e.Cell).Value = $e.Value
}

$LTW_deletes = @("N96")
foreach ($c in $LTW_deletes) {
  $ws.Range($c).ClearContents()
}

Write-Host "Applied $(187) cell edits across 8 sheets"